$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5","D6","D8","D9","D10","D11","D12","D13","D16","D17","D18","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D35","D36","D37","D39","D41","D43","D44","D45","D46","D50")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "73.169.61"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "3.965.21"
$ws.Range("E3").Value = "  -1.85%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "609.64"
$ws.Range("E5").Value = "  +9.01%  "
$ws.Range("D6").Value = "168.27"
$ws.Range("E6").Value = "  +11.72%  "
$ws.Range("E7").Value = "  -2.18%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "0.786"
$ws.Range("E9").Value = "  +2.98%  "
$ws.Range("D10").Value = "0.185"
$ws.Range("E10").Value = "  +7.16%  "
$ws.Range("D11").Value = "55.94"
$ws.Range("E11").Value = "  +3.75%  "
$ws.Range("D12").Value = "0.0000335"
$ws.Range("E12").Value = "  +1.56%  "
$ws.Range("D13").Value = "11.34"
$ws.Range("E13").Value = "  +3.01%  "
$ws.Range("D14").Value = "4.606.82"
$ws.Range("E14").Value = "  -1.89%  "
$ws.Range("D15").Value = "3.980.93"
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("D16").Value = "14.20"
$ws.Range("E16").Value = "  -2.10%  "
$ws.Range("D17").Value = "1.23"
$ws.Range("E17").Value = "  +2.28%  "
$ws.Range("D18").Value = "20.45"
$ws.Range("E18").Value = "  -1.70%  "
$ws.Range("D19").Value = "73.136.78"
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("E20").Value = "  -1.14%  "
$ws.Range("D21").Value = "438.74"
$ws.Range("E21").Value = "  -1.50%  "
$ws.Range("D22").Value = "4.83"
$ws.Range("E22").Value = "  +8.93%  "
$ws.Range("D23").Value = "95.39"
$ws.Range("E23").Value = "  -2.90%  "
$ws.Range("D24").Value = "3.39"
$ws.Range("E24").Value = "  -4.30%  "
$ws.Range("D25").Value = "14.22"
$ws.Range("E25").Value = "  -3.79%  "
$ws.Range("D26").Value = "4.17"
$ws.Range("E26").Value = "  -4.62%  "
$ws.Range("D27").Value = "11.08"
$ws.Range("E27").Value = "  -2.30%  "
$ws.Range("D28").Value = "5.97"
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("D29").Value = "10.51"
$ws.Range("E29").Value = "  -4.40%  "
$ws.Range("D30").Value = "36.05"
$ws.Range("E30").Value = "  -3.22%  "
$ws.Range("D31").Value = "7.81"
$ws.Range("E31").Value = "  -1.62%  "
$ws.Range("D32").Value = "13.85"
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("D33").Value = "0.0000105"
$ws.Range("E33").Value = "  +13.36%  "
$ws.Range("E34").Value = "  -3.68%  "
$ws.Range("D35").Value = "48.21"
$ws.Range("E35").Value = "  -1.52%  "
$ws.Range("D36").Value = "70.40"
$ws.Range("E36").Value = "  +4.37%  "
$ws.Range("D37").Value = "647.34"
$ws.Range("E37").Value = "  -5.60%  "
$ws.Range("E38").Value = "  -4.42%  "
$ws.Range("D39").Value = "3.42"
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("D41").Value = "0.145"
$ws.Range("E41").Value = "  -3.05%  "
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("D43").Value = "0.0485"
$ws.Range("E43").Value = "  -2.63%  "
$ws.Range("D44").Value = "3.18"
$ws.Range("E44").Value = "  -5.77%  "
$ws.Range("D45").Value = "10.48"
$ws.Range("E45").Value = "  -5.06%  "
$ws.Range("D46").Value = "3.10"
$ws.Range("E46").Value = "  +31.55%  "
$ws.Range("E47").Value = "  -2.28%  "
$ws.Range("E48").Value = "  +6.58%  "
$ws.Range("E49").Value = "  +2.79%  "
$ws.Range("D50").Value = "2.58"
$ws.Range("E50").Value = "  -4.46%  "
$ws.Range("E51").Value = "  -4.65%  "
